$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "unit" column (C) for Seedyield (row 2) and Seedyield_bio (row 3)
# to clarify both seed yield traits are standardized to 100% dry matter.
$ws.Range("C2").Value = "dt/ha @100% dry mass"
$ws.Range("C3").Value = "g/m^2 @100% dry mass"
